{"js": "// Replace the trailing empty paragraph (after \"Muaz\") with two narrative\n// paragraphs describing the project-management workflow, plus a new\n// trailing empty paragraph. All three new paragraphs get a first-line\n// indent of 720 twips (36 pt).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph whose text is \"Muaz\" (last bullet of \"Business Plan\n// & Outreach\"), then its immediate successor \u2014 the empty trailing\n// paragraph that the diff replaces.\nlet muazIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Muaz\") {\n    muazIndex = i;\n    break;\n  }\n}\nconst targetPara = paragraphs.items[muazIndex + 1];\n\n// Insert the first new paragraph before the empty trailing paragraph.\nconst para1Text =\n  \"Till the submission of Preliminary Design Review, Project workflow was managed using \" +\n  \"Microsoft Project and all the literature, files, written materials & proposed CAD models were placed in common Google Drive folder.\";\nconst para1 = targetPara.insertParagraph(para1Text, Word.InsertLocation.before);\npara1.firstLineIndent = 36;\n\n// Insert the second new paragraph, still before the (now last) empty paragraph.\nconst para2Text =\n  \"For this Critical Design of UAV and stages after this, \" +\n  \"project is managed using Project Tracker on Google Drive. All files in project including CAD models, reports, simulations, source codes, etc. are managed on GitHub repository. All groups in the team work on different branches in repo.\";\nconst para2 = targetPara.insertParagraph(para2Text, Word.InsertLocation.before);\npara2.firstLineIndent = 36;\n\n// The original trailing empty paragraph remains as the new trailing\n// paragraph; give it the same first-line indent.\ntargetPara.firstLineIndent = 36;\n\nawait context.sync();\n", "ps1": "# Replace the trailing empty paragraph (after \"Muaz\") with two narrative\n# paragraphs describing the project-management workflow, plus a new\n# trailing empty paragraph. All three new paragraphs get a first-line\n# indent of 720 twips (36 pt).\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph whose text is \"Muaz\" (last bullet of \"Business Plan\n# & Outreach\"), then its immediate successor \u2014 the empty trailing\n# paragraph that the diff replaces.\n$muazPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Muaz\") {\n        $muazPara = $p\n        break\n    }\n}\n$targetPara = $muazPara.Next()\n$targetRange = $targetPara.Range\n$targetIndex = $targetPara.Index\n\n$text1 = \"Till the submission of Preliminary Design Review, Project workflow was managed using Microsoft Project and all the literature, files, written materials & proposed CAD models were placed in common Google Drive folder.`r\"\n$text2 = \"For this Critical Design of UAV and stages after this, project is managed using Project Tracker on Google Drive. All files in project including CAD models, reports, simulations, source codes, etc. are managed on GitHub repository. All groups in the team work on different branches in repo.`r\"\n\n# Insert both new paragraphs before the existing trailing empty paragraph.\n$targetRange.InsertBefore($text1 + $text2)\n\n# The two new paragraphs now occupy the original index and the one after\n# it; the original empty paragraph is pushed two slots further.\n$newPara1 = $d.Paragraphs.Item($targetIndex)\n$newPara2 = $d.Paragraphs.Item($targetIndex + 1)\n$newPara3 = $d.Paragraphs.Item($targetIndex + 2)\n\n$newPara1.Format.FirstLineIndent = 36\n$newPara2.Format.FirstLineIndent = 36\n$newPara3.Format.FirstLineIndent = 36\n"}
